# Weekly update: insert this week's price record for
# "Terminal La Palmera de La Serena - Camote" as the new first data row
# (row 4, right after the two most-recent existing entries), pushing the
# older rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 4:24 down to 5:25, inserting a fresh row 4.
$ws.Rows(4).Insert()

# Populate the new weekly record in row 4.
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 45230
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100114002
$ws.Range("G4").Value = "Camote"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 360
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 17000
$ws.Range("M4").Value = 16500
$ws.Range("N4").Value = "`$/malla 18 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 917
$ws.Range("Q4").Value = 18
$ws.Range("R4").Value = "Hortaliza"
